# Auto-generated script applying the cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-CellText $ws 'D2' '28.641.03'
Set-CellText $ws 'E2' '  +0.24%  '
Set-CellText $ws 'D3' '1.874.12'
Set-CellText $ws 'E3' '  +0.62%  '
Set-CellText $ws 'D4' '1.005'
Set-CellText $ws 'E4' '  -0.03%  '
Set-CellText $ws 'D5' '324.99'
Set-CellText $ws 'E5' '  -0.33%  '
Set-CellText $ws 'E6' '  +0.03%  '
Set-CellText $ws 'E7' '  -1.71%  '
Set-CellText $ws 'D8' '0.3839'
Set-CellText $ws 'E8' '  -1.57%  '
Set-CellText $ws 'B9' 'OKB'
Set-CellText $ws 'C9' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-CellText $ws 'D9' '46.64'
Set-CellText $ws 'E9' '  +0.64%  '
Set-CellText $ws 'B10' 'Dogecoin'
Set-CellText $ws 'C10' 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-CellText $ws 'D10' '0.07817'
Set-CellText $ws 'E10' '  -0.87%  '
Set-CellText $ws 'B11' 'Polygon'
Set-CellText $ws 'C11' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-CellText $ws 'D11' '0.9895'
Set-CellText $ws 'E11' '  +2.34%  '
Set-CellText $ws 'B12' 'Solana'
Set-CellText $ws 'C12' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-CellText $ws 'D12' '21.48'
Set-CellText $ws 'E12' '  -3.43%  '
Set-CellText $ws 'B13' 'WrappedEther'
Set-CellText $ws 'C13' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-CellText $ws 'D13' '1.895.96'
Set-CellText $ws 'E13' '  -2.83%  '
Set-CellText $ws 'B14' 'Chainlink'
Set-CellText $ws 'C14' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-CellText $ws 'D14' '6.907'
Set-CellText $ws 'E14' '  -0.16%  '
Set-CellText $ws 'B15' 'Polkadot'
Set-CellText $ws 'C15' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-CellText $ws 'D15' '5.637'
Set-CellText $ws 'E15' '  -1.26%  '
Set-CellText $ws 'B16' 'TRON'
Set-CellText $ws 'C16' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-CellText $ws 'D16' '0.06940'
Set-CellText $ws 'E16' '  -0.34%  '
Set-CellText $ws 'B17' 'Litecoin'
Set-CellText $ws 'C17' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-CellText $ws 'D17' '86.62'
Set-CellText $ws 'E17' '  -1.73%  '
Set-CellText $ws 'B18' 'BinanceUSD'
Set-CellText $ws 'C18' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-CellText $ws 'D18' '1.006'
Set-CellText $ws 'E18' '  -0.04%  '
Set-CellText $ws 'B19' 'ShibaInu'
Set-CellText $ws 'C19' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-CellText $ws 'D19' '0.000009950'
Set-CellText $ws 'E19' '  -0.86%  '
Set-CellText $ws 'B20' 'Avalanche'
Set-CellText $ws 'C20' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-CellText $ws 'D20' '16.72'
Set-CellText $ws 'E20' '  -1.28%  '
Set-CellText $ws 'B21' 'Dai'
Set-CellText $ws 'C21' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-CellText $ws 'D21' '1.004'
Set-CellText $ws 'E21' '  -0.03%  '
Set-CellText $ws 'B22' 'WrappedBTC'
Set-CellText $ws 'C22' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-CellText $ws 'D22' '28.624.08'
Set-CellText $ws 'E22' '  +0.01%  '
Set-CellText $ws 'B23' 'Uniswap'
Set-CellText $ws 'C23' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-CellText $ws 'D23' '5.258'
Set-CellText $ws 'E23' '  -0.82%  '
Set-CellText $ws 'B24' 'Cosmos'
Set-CellText $ws 'C24' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-CellText $ws 'D24' '10.88'
Set-CellText $ws 'E24' '  -1.51%  '
Set-CellText $ws 'B25' 'Toncoin'
Set-CellText $ws 'C25' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-CellText $ws 'D25' '2.064'
Set-CellText $ws 'E25' '  -2.77%  '
Set-CellText $ws 'B26' 'WrappedliquidstakedEther2.0'
Set-CellText $ws 'C26' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-CellText $ws 'D26' '2.096.03'
Set-CellText $ws 'E26' '  -2.87%  '
Set-CellText $ws 'B27' 'Monero'
Set-CellText $ws 'C27' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-CellText $ws 'D27' '154.46'
Set-CellText $ws 'E27' '  +0.49%  '
Set-CellText $ws 'B28' 'EthereumClassic'
Set-CellText $ws 'C28' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-CellText $ws 'D28' '19.12'
Set-CellText $ws 'E28' '  -0.91%  '
Set-CellText $ws 'B29' 'InternetComputer(DFINITY)'
Set-CellText $ws 'C29' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-CellText $ws 'D29' '5.653'
Set-CellText $ws 'E29' '  -0.75%  '
Set-CellText $ws 'B30' 'BitcoinCash'
Set-CellText $ws 'C30' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-CellText $ws 'D30' '117.54'
Set-CellText $ws 'E30' '  -1.27%  '
Set-CellText $ws 'B31' 'LidoDAOToken'
Set-CellText $ws 'C31' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-CellText $ws 'D31' '1.877'
Set-CellText $ws 'E31' '  -5.69%  '
Set-CellText $ws 'B32' 'Stellar'
Set-CellText $ws 'C32' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-CellText $ws 'D32' '0.09261'
Set-CellText $ws 'E32' '  -0.85%  '
Set-CellText $ws 'B33' 'ImmutableX'
Set-CellText $ws 'C33' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-CellText $ws 'D33' '0.9022'
Set-CellText $ws 'E33' '  -2.87%  '
Set-CellText $ws 'B34' 'Filecoin'
Set-CellText $ws 'C34' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-CellText $ws 'D34' '5.269'
Set-CellText $ws 'E34' '  -0.83%  '
Set-CellText $ws 'B35' 'ARBITRUM'
Set-CellText $ws 'C35' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-CellText $ws 'D35' '1.319'
Set-CellText $ws 'E35' '  -1.38%  '
Set-CellText $ws 'B36' 'HuobiToken'
Set-CellText $ws 'C36' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-CellText $ws 'D36' '3.253'
Set-CellText $ws 'E36' '  -3.00%  '
Set-CellText $ws 'B37' 'Hedera'
Set-CellText $ws 'C37' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-CellText $ws 'D37' '0.05654'
Set-CellText $ws 'E37' '  -2.80%  '
Set-CellText $ws 'B38' 'TrustWalletToken'
Set-CellText $ws 'C38' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-CellText $ws 'D38' '1.149'
Set-CellText $ws 'E38' '  -0.09%  '
Set-CellText $ws 'B39' 'VeChain'
Set-CellText $ws 'C39' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-CellText $ws 'D39' '0.02041'
Set-CellText $ws 'E39' '  -3.67%  '
Set-CellText $ws 'B40' 'FraxShare'
Set-CellText $ws 'C40' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-CellText $ws 'D40' '7.604'
Set-CellText $ws 'E40' '  -3.43%  '
Set-CellText $ws 'B41' 'TheSandbox'
Set-CellText $ws 'C41' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-CellText $ws 'D41' '0.5550'
Set-CellText $ws 'E41' '  -1.55%  '
Set-CellText $ws 'B42' 'Algorand'
Set-CellText $ws 'C42' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-CellText $ws 'D42' '0.1764'
Set-CellText $ws 'E42' '  -0.70%  '
Set-CellText $ws 'B43' 'Aptos'
Set-CellText $ws 'C43' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-CellText $ws 'D43' '9.577'
Set-CellText $ws 'E43' '  -3.22%  '
Set-CellText $ws 'B44' 'Cronos'
Set-CellText $ws 'C44' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-CellText $ws 'D44' '0.07153'
Set-CellText $ws 'E44' '  -1.09%  '
Set-CellText $ws 'D45' '11.52'
Set-CellText $ws 'E45' '  -1.40%  '
Set-CellText $ws 'B46' 'Decentraland'
Set-CellText $ws 'C46' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-CellText $ws 'D46' '0.5237'
Set-CellText $ws 'E46' '  -1.29%  '
Set-CellText $ws 'B47' 'RenderToken'
Set-CellText $ws 'C47' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-CellText $ws 'D47' '2.153'
Set-CellText $ws 'E47' '  -0.06%  '
Set-CellText $ws 'B48' 'WEMIXToken'
Set-CellText $ws 'C48' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-CellText $ws 'D48' '1.114'
Set-CellText $ws 'E48' '  -1.99%  '
Set-CellText $ws 'B49' 'NEARProtocol'
Set-CellText $ws 'C49' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-CellText $ws 'D49' '1.798'
Set-CellText $ws 'E49' '  -2.27%  '
Set-CellText $ws 'B50' 'Quant'
Set-CellText $ws 'C50' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-CellText $ws 'D50' '111.20'
Set-CellText $ws 'E50' '  -1.70%  '
Set-CellText $ws 'B51' 'MXToken'
Set-CellText $ws 'C51' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-CellText $ws 'D51' '2.434'
Set-CellText $ws 'E51' '  +3.84%  '
